# -----------------------------------------------------------------------
# Scheduled runner update: refresh currentAveragePrice* / Leve profit
# columns (H:N) on each job sheet of the Phoenix_Profits workbook with
# newly pulled market-board data.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Cells.Item(18, 10).Value = 3009.8  # J18
$ws.Cells.Item(18, 12).Value = 3009.8  # L18
$ws.Cells.Item(18, 8).Value = 2023.0714  # H18
$ws.Cells.Item(18, 11).Value = 1474.8889  # K18
$ws.Cells.Item(18, 9).Value = 1474.8889  # I18
$ws.Cells.Item(18, 13).Value = -1190.8889  # M18
$ws.Cells.Item(18, 14).Value = -3577.8  # N18
# Row 40
$ws.Cells.Item(40, 11).Value = 4545.4546  # K40
$ws.Cells.Item(40, 12).Value = 4700  # L40
$ws.Cells.Item(40, 13).Value = -4370.4546  # M40
$ws.Cells.Item(40, 10).Value = 4700  # J40
$ws.Cells.Item(40, 9).Value = 4545.4546  # I40
$ws.Cells.Item(40, 8).Value = 4593.75  # H40
$ws.Cells.Item(40, 14).Value = -5050  # N40
# Row 62
$ws.Cells.Item(62, 9).Value = 4945.75  # I62
$ws.Cells.Item(62, 8).Value = 5597.5  # H62
$ws.Cells.Item(62, 14).Value = -7497.25  # N62
$ws.Cells.Item(62, 10).Value = 6249.25  # J62
$ws.Cells.Item(62, 12).Value = 6249.25  # L62
$ws.Cells.Item(62, 13).Value = -4321.75  # M62
$ws.Cells.Item(62, 11).Value = 4945.75  # K62
# Row 65
$ws.Cells.Item(65, 12).Value = 31246.25  # L65
$ws.Cells.Item(65, 8).Value = 5597.5  # H65
$ws.Cells.Item(65, 11).Value = 24728.75  # K65
$ws.Cells.Item(65, 10).Value = 6249.25  # J65
$ws.Cells.Item(65, 14).Value = -37486.25  # N65
$ws.Cells.Item(65, 9).Value = 4945.75  # I65
$ws.Cells.Item(65, 13).Value = -21608.75  # M65
# Row 135
$ws.Cells.Item(135, 13).Value = -51070.00169999999  # M135
$ws.Cells.Item(135, 8).Value = 6687.5  # H135
$ws.Cells.Item(135, 11).Value = 53605.00169999999  # K135
$ws.Cells.Item(135, 9).Value = 5956.1113  # I135
# Row 137
$ws.Cells.Item(137, 8).Value = 2712.7812  # H137
$ws.Cells.Item(137, 13).Value = -2885.700000000001  # M137
$ws.Cells.Item(137, 11).Value = 5435.700000000001  # K137
$ws.Cells.Item(137, 9).Value = 1811.9  # I137

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 13).Value = -23258.596  # M32
$ws.Cells.Item(32, 9).Value = 23545.596  # I32
$ws.Cells.Item(32, 11).Value = 23545.596  # K32
$ws.Cells.Item(32, 8).Value = 23457.049  # H32
# Row 61
$ws.Cells.Item(61, 13).Value = -3097.9473  # M61
$ws.Cells.Item(61, 9).Value = 3309.9473  # I61
$ws.Cells.Item(61, 8).Value = 30157.195  # H61
$ws.Cells.Item(61, 11).Value = 3309.9473  # K61
# Row 122
$ws.Cells.Item(122, 13).Value = -5086.428400000001  # M122
$ws.Cells.Item(122, 11).Value = 7536.428400000001  # K122
$ws.Cells.Item(122, 9).Value = 2512.1428  # I122
$ws.Cells.Item(122, 8).Value = 38598.08  # H122
# Row 136
$ws.Cells.Item(136, 11).Value = 9929.841899999999  # K136
$ws.Cells.Item(136, 8).Value = 30157.195  # H136
$ws.Cells.Item(136, 9).Value = 3309.9473  # I136
$ws.Cells.Item(136, 13).Value = -7379.841899999999  # M136

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Cells.Item(64, 9).Value = 387.2  # I64
$ws.Cells.Item(64, 8).Value = 1401.4736  # H64
$ws.Cells.Item(64, 13).Value = -162.2  # M64
$ws.Cells.Item(64, 11).Value = 387.2  # K64
# Row 67
$ws.Cells.Item(67, 11).Value = 387.2  # K67
$ws.Cells.Item(67, 8).Value = 1401.4736  # H67
$ws.Cells.Item(67, 13).Value = 392.8  # M67
$ws.Cells.Item(67, 9).Value = 387.2  # I67
# Row 86
$ws.Cells.Item(86, 9).Value = 869.0714  # I86
$ws.Cells.Item(86, 11).Value = 869.0714  # K86
$ws.Cells.Item(86, 13).Value = 253.9286  # M86
$ws.Cells.Item(86, 8).Value = 45893.35  # H86
# Row 89
$ws.Cells.Item(89, 9).Value = 869.0714  # I89
$ws.Cells.Item(89, 11).Value = 4345.357  # K89
$ws.Cells.Item(89, 8).Value = 45893.35  # H89
$ws.Cells.Item(89, 13).Value = 1270.643  # M89

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 9).Value = 2222.9583  # I31
$ws.Cells.Item(31, 11).Value = 2222.9583  # K31
$ws.Cells.Item(31, 13).Value = -1927.9583  # M31
$ws.Cells.Item(31, 8).Value = 2595.3416  # H31
# Row 34
$ws.Cells.Item(34, 9).Value = 2222.9583  # I34
$ws.Cells.Item(34, 8).Value = 2595.3416  # H34
$ws.Cells.Item(34, 13).Value = -2020.9583  # M34
$ws.Cells.Item(34, 11).Value = 2222.9583  # K34
# Row 58
$ws.Cells.Item(58, 9).Value = 2870.8333  # I58
$ws.Cells.Item(58, 11).Value = 2870.8333  # K58
$ws.Cells.Item(58, 10).Value = 7714.6665  # J58
$ws.Cells.Item(58, 14).Value = -8120.6665  # N58
$ws.Cells.Item(58, 13).Value = -2667.8333  # M58
$ws.Cells.Item(58, 12).Value = 7714.6665  # L58
$ws.Cells.Item(58, 8).Value = 3839.6  # H58
# Row 136
$ws.Cells.Item(136, 10).Value = 7714.6665  # J136
$ws.Cells.Item(136, 14).Value = -28243.9995  # N136
$ws.Cells.Item(136, 12).Value = 23143.9995  # L136
$ws.Cells.Item(136, 11).Value = 8612.499899999999  # K136
$ws.Cells.Item(136, 8).Value = 3839.6  # H136
$ws.Cells.Item(136, 9).Value = 2870.8333  # I136
$ws.Cells.Item(136, 13).Value = -6062.499899999999  # M136

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Cells.Item(62, 9).Value = 8200  # I62
$ws.Cells.Item(62, 8).Value = 9485  # H62
$ws.Cells.Item(62, 14).Value = -31369  # N62
$ws.Cells.Item(62, 10).Value = 9999  # J62
$ws.Cells.Item(62, 12).Value = 29997  # L62
$ws.Cells.Item(62, 13).Value = -23914  # M62
$ws.Cells.Item(62, 11).Value = 24600  # K62
# Row 65
$ws.Cells.Item(65, 12).Value = 89991  # L65
$ws.Cells.Item(65, 8).Value = 9485  # H65
$ws.Cells.Item(65, 11).Value = 73800  # K65
$ws.Cells.Item(65, 10).Value = 9999  # J65
$ws.Cells.Item(65, 14).Value = -96855  # N65
$ws.Cells.Item(65, 9).Value = 8200  # I65
$ws.Cells.Item(65, 13).Value = -70368  # M65
# Row 74
$ws.Cells.Item(74, 8).Value = 10000  # H74
$ws.Cells.Item(74, 9).Value = 10000  # I74
$ws.Cells.Item(74, 13).Value = -28939  # M74
$ws.Cells.Item(74, 11).Value = 30000  # K74
# Row 77
$ws.Cells.Item(77, 13).Value = -84696  # M77
$ws.Cells.Item(77, 8).Value = 10000  # H77
$ws.Cells.Item(77, 9).Value = 10000  # I77
$ws.Cells.Item(77, 11).Value = 90000  # K77
# Row 80
$ws.Cells.Item(80, 8).Value = 16664  # H80
$ws.Cells.Item(80, 10).Value = 15000  # J80
$ws.Cells.Item(80, 14).Value = -46872  # N80
$ws.Cells.Item(80, 12).Value = 45000  # L80
# Row 81
$ws.Cells.Item(81, 14).ClearContents() | Out-Null  # N81: remove (was -8240)
$ws.Cells.Item(81, 11).Value = 3991.0002  # K81
$ws.Cells.Item(81, 9).Value = 1330.3334  # I81
$ws.Cells.Item(81, 8).Value = 1330.3334  # H81
$ws.Cells.Item(81, 13).Value = -2868.0002  # M81
$ws.Cells.Item(81, 12).Value = 0  # L81
$ws.Cells.Item(81, 10).Value = 0  # J81
# Row 83
$ws.Cells.Item(83, 8).Value = 16664  # H83
$ws.Cells.Item(83, 12).Value = 135000  # L83
$ws.Cells.Item(83, 14).Value = -144360  # N83
$ws.Cells.Item(83, 10).Value = 15000  # J83
# Row 84
$ws.Cells.Item(84, 13).Value = -6357.000599999999  # M84
$ws.Cells.Item(84, 14).ClearContents() | Out-Null  # N84: remove (was -29214)
$ws.Cells.Item(84, 10).Value = 0  # J84
$ws.Cells.Item(84, 11).Value = 11973.0006  # K84
$ws.Cells.Item(84, 12).Value = 0  # L84
$ws.Cells.Item(84, 8).Value = 1330.3334  # H84
$ws.Cells.Item(84, 9).Value = 1330.3334  # I84
# Row 112
$ws.Cells.Item(112, 8).Value = 9341.333000000001  # H112
$ws.Cells.Item(112, 9).Value = 9341.333000000001  # I112
$ws.Cells.Item(112, 11).Value = 28023.999  # K112
$ws.Cells.Item(112, 13).Value = -26915.999  # M112

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 11).Value = 75.17646999999999  # K2
$ws.Cells.Item(2, 12).Value = 48.5  # L2
$ws.Cells.Item(2, 13).Value = 37.82353000000001  # M2
$ws.Cells.Item(2, 8).Value = 66.64  # H2
$ws.Cells.Item(2, 10).Value = 48.5  # J2
$ws.Cells.Item(2, 9).Value = 75.17646999999999  # I2
$ws.Cells.Item(2, 14).Value = -274.5  # N2
# Row 80
$ws.Cells.Item(80, 8).Value = 23826.834  # H80
$ws.Cells.Item(80, 10).Value = 28370  # J80
$ws.Cells.Item(80, 14).Value = -30366  # N80
$ws.Cells.Item(80, 12).Value = 28370  # L80
# Row 83
$ws.Cells.Item(83, 8).Value = 23826.834  # H83
$ws.Cells.Item(83, 12).Value = 141850  # L83
$ws.Cells.Item(83, 14).Value = -151834  # N83
$ws.Cells.Item(83, 10).Value = 28370  # J83
# Row 122
$ws.Cells.Item(122, 12).Value = 12108  # L122
$ws.Cells.Item(122, 14).Value = -17008  # N122
$ws.Cells.Item(122, 10).Value = 4036  # J122
$ws.Cells.Item(122, 8).Value = 2545.1667  # H122
# Row 132
$ws.Cells.Item(132, 11).Value = 11784.3075  # K132
$ws.Cells.Item(132, 8).Value = 3942  # H132
$ws.Cells.Item(132, 13).Value = -9254.307499999999  # M132
$ws.Cells.Item(132, 9).Value = 3928.1025  # I132

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 2079.423  # H16
$ws.Cells.Item(16, 9).Value = 1046.3043  # I16
$ws.Cells.Item(16, 10).Value = 10000  # J16
$ws.Cells.Item(16, 11).Value = 1046.3043  # K16
$ws.Cells.Item(16, 13).Value = -876.3043  # M16
$ws.Cells.Item(16, 14).Value = -10340  # N16
$ws.Cells.Item(16, 12).Value = 10000  # L16
# Row 46
$ws.Cells.Item(46, 8).Value = 3873.587  # H46
$ws.Cells.Item(46, 9).Value = 4094.5  # I46
$ws.Cells.Item(46, 12).Value = 3863.5454  # L46
$ws.Cells.Item(46, 13).Value = -3906.5  # M46
$ws.Cells.Item(46, 10).Value = 3863.5454  # J46
$ws.Cells.Item(46, 11).Value = 4094.5  # K46
$ws.Cells.Item(46, 14).Value = -4239.5454  # N46
# Row 61
$ws.Cells.Item(61, 13).Value = -7152  # M61
$ws.Cells.Item(61, 9).Value = 7354  # I61
$ws.Cells.Item(61, 8).Value = 7828.3335  # H61
$ws.Cells.Item(61, 11).Value = 7354  # K61
# Row 93
$ws.Cells.Item(93, 8).Value = 2167.9443  # H93
$ws.Cells.Item(93, 9).Value = 1849.6666  # I93
$ws.Cells.Item(93, 13).Value = -601.6666  # M93
$ws.Cells.Item(93, 11).Value = 1849.6666  # K93
# Row 113
$ws.Cells.Item(113, 8).Value = 7828.3335  # H113
$ws.Cells.Item(113, 11).Value = 7354  # K113
$ws.Cells.Item(113, 13).Value = -5184  # M113
$ws.Cells.Item(113, 9).Value = 7354  # I113
# Row 132
$ws.Cells.Item(132, 11).Value = 8120.750100000001  # K132
$ws.Cells.Item(132, 8).Value = 3376.3547  # H132
$ws.Cells.Item(132, 13).Value = -5590.750100000001  # M132
$ws.Cells.Item(132, 9).Value = 2706.9167  # I132

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 311.57144  # H107
$ws.Cells.Item(107, 11).Value = 934.71432  # K107
$ws.Cells.Item(107, 9).Value = 311.57144  # I107
$ws.Cells.Item(107, 13).Value = 985.28568  # M107

Write-Host "Applied 195 cell updates across 8 sheets."
